$p = $ppt.ActivePresentation

# Move slide 10 ("Bias in artificial intelligence") to the end of the deck
# (after the current last slide, position 12).
$count = $p.Slides.Count
$p.Slides.Item(10).MoveTo($count)
